# contact.xlsx — "Add files via upload"
#
# The sheet previously held a header row plus one sample data row
# (ID 741710024 / Last anketa 0.5). The re-uploaded version only keeps
# the header row: the data row is removed entirely and the running
# counter in Z1 drops from 3 to 2. The view is also left scrolled over
# to column B, zoomed to 70%, with Z9 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the second row (A2="741710024", B2="0.5") completely, shifting
# sheetData back down to just the header row (dimension becomes A1:Z1).
$ws.Rows("2:2").Delete() | Out-Null

# The trailing counter cell goes from 3 to 2.
$ws.Range("Z1").Value = 2

# Match the saved view: scrolled so column B is left-most, zoomed to
# 70%, with Z9 as the active/selected cell.
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.Zoom = 70
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("Z9").Select() | Out-Null
